# Rename the sheet "SRAB" to "RNAME" and make it the active/selected sheet
# (moving the active tab away from "CCC").

$wb = $excel.ActiveWorkbook

$sheet = $wb.Worksheets.Item("SRAB")
$sheet.Name = "RNAME"

$sheet.Activate()
$sheet.Select()
